# Update countries & provincias Spain
# Applies the refreshed COVID-19 dataset numbers and swaps the
# Birmania / Jordania ranking (Birmania's updated total now edges
# ahead of Jordania's), plus bumps the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 05:59"

# --- Pakistan (row 20) ---
$ws.Range("B20").Value = 303634
$ws.Range("C20").Value = 545
$ws.Range("D20").Value = 291169
$ws.Range("E20").Value = 6066
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 6399

# --- Belgica (row 39) ---
$ws.Range("B39").Value = 95948
$ws.Range("C39").Value = 1153
$ws.Range("D39").Value = 18810
$ws.Range("E39").Value = 67203
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 5
$ws.Range("H39").Value = 9935

# --- Honduras (row 50) ---
$ws.Range("B50").Value = 69660
$ws.Range("C50").Value = 1040
$ws.Range("D50").Value = 19983
$ws.Range("E50").Value = 47575
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 15
$ws.Range("H50").Value = 2102

# --- Venezuela (row 53) ---
$ws.Range("B53").Value = 63416
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 51274
$ws.Range("E53").Value = 11631
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 511

# --- Birmania / Jordania swap (rows 127-128) ---
# Birmania's refreshed numbers now outrank Jordania's, so Birmania
# moves up into row 127 (with new data) and Jordania drops to row 128
# (keeping the figures Birmania's old row 127 used to hold).
$ws.Range("A127").Value = "Birmania"
$ws.Range("B127").Value = 3894
$ws.Range("C127").Value = 73
$ws.Range("D127").Value = 908
$ws.Range("E127").Value = 2940
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 6
$ws.Range("H127").Value = 46

$ws.Range("A128").Value = "Jordania"
$ws.Range("B128").Value = 3852
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 2349
$ws.Range("E128").Value = 1477
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 26

# --- San Martin (Parte Holandesa) (row 173) ---
$ws.Range("B173").Value = 557
$ws.Range("C173").Value = 8
$ws.Range("D173").Value = 477
$ws.Range("E173").Value = 61
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 19
